# Applies the edit described in the commit:
#  - Update the cached "today" date field shown in the Date placeholder on
#    the slide master and every slide layout, from 4/16/2018 to 4/20/2018.
#  - On slide 15 ("DDoS Attack: On a Web Server hosted on a Virtual
#    Machine"), relabel the two comparison callouts:
#      "DPS : Basic"    -> "Without DPS"
#      "DPS : Standard" -> "With DPS Standard" (and widen the textbox so the
#      longer label still fits)
#  - Delete the last slide in the deck (the DDoS Protection Service
#    offerings reference slide).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Refresh the cached date-field text (type 16 == ppPlaceholderDate) on
#    the slide master and on every slide layout.
# ---------------------------------------------------------------------
$newDate = "4/20/2018"

function Update-DatePlaceholder($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        $isDate = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) { $isDate = $true }
        } catch {
        }
        if ($isDate -and $sh.HasTextFrame) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Update-DatePlaceholder $layouts.Item($i).Shapes
}

# ---------------------------------------------------------------------
# 2) Relabel the DPS comparison callouts on slide 15.
# ---------------------------------------------------------------------
$slide15 = $p.Slides.Item(15)

for ($k = 1; $k -le $slide15.Shapes.Count; $k++) {
    $sh = $slide15.Shapes.Item($k)
    if (-not $sh.HasTextFrame) { continue }
    $txt = $sh.TextFrame.TextRange.Text
    if ($txt -eq "DPS : Basic") {
        $sh.TextFrame.TextRange.Text = "Without DPS"
    } elseif ($txt -eq "DPS : Standard") {
        $sh.TextFrame.TextRange.Text = "With DPS Standard"
        # Widen the textbox (1638955 EMU -> 2010837 EMU) to fit the longer
        # label; Shape.Width is in points (EMU / 12700).
        $sh.Width = 158.333626
    }
}

# ---------------------------------------------------------------------
# 3) Drop the last slide from the deck.
# ---------------------------------------------------------------------
$p.Slides.Item($p.Slides.Count).Delete()
